$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-unused trailing columns (I:N) ---------------------------
# Old table went out to column N (2008-2018 + 11 data points). The new
# table only needs columns D:H (2017-2020 + "2021*"), so blow away I:N
# entirely (shifts the used range/dimension back down to column H).
$ws.Range("I1:N8").EntireColumn.Delete()

# --- Row 1: title (drop the trailing "*") ----------------------------------
$ws.Range("A1").Value = "9.5.1 ИДП га болгон тажрыйбалык-конструктордук жумуштун жана илимий изилдөөнүн чыгымдарынын үлүшү"
$ws.Range("B1").Value = "9.5.1  Доля расходов на научно-исследовательские и опытно-конструкторские работы в ВВП"
$ws.Range("C1").Value = "9.5.1 Research and development expenditure as a proportion of GDP"
$ws.Rows("1:1").RowHeight = 43.5

# --- Row 4: header labels + years -------------------------------------------
$ws.Range("D4").Value = 2017
$ws.Range("E4").Value = 2018
$ws.Range("F4").Value = 2019
$ws.Range("G4").Value = 2020
$ws.Range("H4").Value = "2021*"

# Copy the year-cell style (D4) onto F4:H4 (F4/G4 used to be plain s=2
# empties, H4 is the brand-new "2021*" footnote-year cell) and then just
# right-align H4's text (numbers default to right-aligned, "2021*" is a
# string so it needs the alignment spelled out explicitly).
$ws.Range("D4").Copy()
$ws.Range("F4:H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H4").HorizontalAlignment = -4152

# --- Row 5: data values -------------------------------------------------
$ws.Range("D5").Value = 0.11
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 0.09
$ws.Range("G5").Value = 0.09
$ws.Range("H5").Value = 0.08

# --- Row 6: footnote (now "*preliminary data" instead of "*per MF KR data") -
$ws.Range("A6").Value = "*алдын алаа маалыматтар"
$ws.Range("B6").Value = "*предварительные данные"
$ws.Range("C6").Value = "*preliminary data"

# --- Remove the leftover stray cells in columns F:H on the short rows -----
# (rows 1,2,6,7,8 only span A:E in the new layout; rows 3,4,5 keep A:H)
$ws.Range("F1:H2").Clear()
$ws.Range("F6:H8").Clear()

$wb.Save()
